$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rows 2-6 (LOW* group): C -0.5 -> -10, D -0.99 -> -20
$ws.Range("C2:C6").Value = -10
$ws.Range("D2:D6").Value = -20

# Rows 7-11 (MED* group): B 0.5 -> 10, D -0.99 -> -20
$ws.Range("B7:B11").Value = 10
$ws.Range("D7:D11").Value = -20

# Rows 12-16 (HIGH* group): B 0.99 -> 20, C 0.5 -> 10, D -0.99 -> -20
$ws.Range("B12:B16").Value = 20
$ws.Range("C12:C16").Value = 10
$ws.Range("D12:D16").Value = -20

# Update the active selection to match the saved view state
$ws.Range("H18").Select() | Out-Null
